# Update "想去人数" (F column) counts for a new scrape snapshot.
# Values scraped anew differ slightly from the previous snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 13960
$ws1.Range("F7").Value = 1041
$ws1.Range("F8").Value = 13911
$ws1.Range("F9").Value = 14879
$ws1.Range("F25").Value = 5788
$ws1.Range("F27").Value = 1063
$ws1.Range("F28").Value = 5431
$ws1.Range("F32").Value = 301

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 13961
$ws4.Range("F8").Value = 1041
$ws4.Range("F9").Value = 13911
$ws4.Range("F10").Value = 14879
$ws4.Range("F26").Value = 5788
$ws4.Range("F28").Value = 1063
$ws4.Range("F29").Value = 5431
$ws4.Range("F33").Value = 301
